$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Haba at Terminal
# Hortofrutícola Agro Chillán. Insert a row above the current row 52 so the
# new (most recent) record lands there, pushing the existing rows 52-69
# down to 53-70.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with the new record's data.
$ws.Cells.Item(52,1).Value = 7
$ws.Cells.Item(52,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(52,3).Value = "Ñuble"
$ws.Cells.Item(52,4).Value = 44900
$ws.Cells.Item(52,5).Value = 16
$ws.Cells.Item(52,6).Value = 100112026
$ws.Cells.Item(52,7).Value = "Haba"
$ws.Cells.Item(52,8).Value = "Sin especificar"
$ws.Cells.Item(52,9).Value = "Primera"
$ws.Cells.Item(52,10).Value = 80
$ws.Cells.Item(52,11).Value = 12000
$ws.Cells.Item(52,12).Value = 12000
$ws.Cells.Item(52,13).Value = 12000
$ws.Cells.Item(52,14).Value = "$/saco 25 kilos"
$ws.Cells.Item(52,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(52,16).Value = 480
$ws.Cells.Item(52,17).Value = 25
$ws.Cells.Item(52,18).Value = "Hortaliza"
